# Commit: Thu, May 28, 2020  5:05:00 AM
#
# Two changes are applied:
#
#  1. Slide 6's table (the graphicFrame holding the "SOURCES OF FINANCE"
#     table) switches its table style from the custom/local style
#     {EF6B999C-D8FE-4FAD-94DB-403141D5CFCA} to the built-in style
#     {F0248284-0730-475D-B70D-199AA28A2564}.
#
#  2. The deck's theme colour palette is swapped from the "Integral"
#     palette to the "Office Theme" palette (the font scheme and format
#     scheme are identical between the two themes in this deck, so only
#     the 12 colour-scheme slots actually need to change).

$p = $ppt.ActivePresentation

# --- 1) Table style on Slide 6 ------------------------------------------
$slide6   = $p.Slides.Item(6)
$tblShape = $slide6.Shapes.Item(2)
$tblShape.Table.ApplyStyle("{F0248284-0730-475D-B70D-199AA28A2564}")

# --- 2) Theme colour scheme: Integral -> Office Theme -------------------
$theme = $p.SlideMaster.Theme
$tcs   = $theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0          # dk1      000000
$tcs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388    # dk2      44546A
$tcs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407      # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308   # accent5  4472C4
$tcs.Colors(10).RGB = 4697456    # accent6  70AD47
$tcs.Colors(11).RGB = 12673797   # hlink    0563C1
$tcs.Colors(12).RGB = 7491477    # folHlink 954F72
